$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the "type" column values for rows 3 and 4 (First Sip Cafe -> "B", Bang Bang Pie & Biscuits -> "C")
$ws.Range("B3").Value = "B"
$ws.Range("B4").Value = "C"

# Add new header cells K1:N1 mirroring the place_id/formatted_address/latitude/longitude headers
$ws.Range("K1").Value = "place_id"
$ws.Range("L1").Value = "formatted_address"
$ws.Range("M1").Value = "latitude"
$ws.Range("N1").Value = "longitude"

# Row 2 - Center on Halsted: duplicate place_id/formatted_address/lat/long into K2:N2
$ws.Range("K2").Value = "ChIJx_VyuLDTD4gR_TjFLHXJdpY"
$ws.Range("L2").Value = "3656 N Halsted St, Chicago, IL 60613, United States"
$ws.Range("M2").Value = 41.9489773
$ws.Range("N2").Value = -87.6497045

# Row 3 - First Sip Cafe: duplicate place_id/formatted_address/lat/long into K3:N3
$ws.Range("K3").Value = "ChIJ33YxJNXTD4gRWe6jz5NdNpo"
$ws.Range("L3").Value = "1057 W Argyle St, Chicago, IL 60640, United States"
$ws.Range("M3").Value = 41.9731758
$ws.Range("N3").Value = -87.6572092

# Row 4 - Bang Bang Pie & Biscuits: duplicate place_id/formatted_address/lat/long into K4:N4
$ws.Range("K4").Value = "ChIJt6mzLZ7SD4gR8bpjii62t_Q"
$ws.Range("L4").Value = "2051 N California Ave, Chicago, IL 60647, United States"
$ws.Range("M4").Value = 41.9190202
$ws.Range("N4").Value = -87.6971201
